# Daily attendance processing
#
# The "Recorded By" column (G) stores a comma-separated list of the
# users/services that recorded a session. For every row whose list is
# exactly "dnasr281@gmail.com, System" or "system, backup@backdoor.com,
# System", rotate the list left by one position (move the first entry to
# the end) to reflect the latest processing pass order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System" -or $val -eq "system, backup@backdoor.com, System") {
        $parts = $val -split ", "
        $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
        $cell.Value = $rotated
        $changed++
    }
}

Write-Host "Recorded By values rotated: $changed"
